$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Update status text from "Ready for handoff" to "Handback transform failed"
$overview.Range("B2").Value = "Handback transform failed"
$overview.Range("C2").Value = "Handback transform failed"

$zhcn.Range("B2").Value = "Handback transform failed"
$dede.Range("B2").Value = "Handback transform failed"

# Add Error Detail text for the handback mismatch in each language sheet
$zhcn.Range("J2").Value = "Handback file name: qswfh3lj.vvt is different with handoff file name: 30118e80-996d-4d7e-9af6-3efa4d887583.57a95008bb74c92f77334feb40d9fb30bba13abe.zh-cn."
$dede.Range("J2").Value = "Handback file name: qswfh3lj.vvt is different with handoff file name: 30118e80-996d-4d7e-9af6-3efa4d887583.57a95008bb74c92f77334feb40d9fb30bba13abe.de-de."
